$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "66.211.39"
Set-TextValue "E2" "  -1.39%  "
Set-TextValue "D3" "3.433.16"
Set-TextValue "E3" "  -0.80%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "582.52"
Set-TextValue "E5" "  -1.82%  "
Set-TextValue "D6" "174.05"
Set-TextValue "E6" "  -3.02%  "
Set-TextValue "D8" "0.591"
Set-TextValue "E8" "  -3.40%  "
Set-TextValue "D9" "3.429.27"
Set-TextValue "E9" "  -0.89%  "
Set-TextValue "E10" "  -6.27%  "
Set-TextValue "D11" "6.87"
Set-TextValue "E11" "  -1.15%  "
Set-TextValue "D12" "0.409"
Set-TextValue "E12" "  -4.89%  "
Set-TextValue "D13" "4.025.88"
Set-TextValue "E13" "  -0.88%  "
Set-TextValue "E14" "  -0.56%  "
Set-TextValue "D15" "29.86"
Set-TextValue "E15" "  -6.32%  "
Set-TextValue "D16" "66.217.26"
Set-TextValue "E16" "  -1.40%  "
Set-TextValue "D17" "0.0000170"
Set-TextValue "E17" "  -3.54%  "
Set-TextValue "D18" "3.439.57"
Set-TextValue "E18" "  -0.62%  "
Set-TextValue "D19" "5.87"
Set-TextValue "E19" "  -5.25%  "
Set-TextValue "D20" "13.70"
Set-TextValue "E20" "  -3.16%  "
Set-TextValue "D21" "373.21"
Set-TextValue "E21" "  -4.90%  "
Set-TextValue "D22" "7.72"
Set-TextValue "E22" "  -2.35%  "
Set-TextValue "D23" "0.997"
Set-TextValue "E23" "  -0.25%  "
Set-TextValue "D24" "71.80"
Set-TextValue "E24" "  +0.25%  "
Set-TextValue "E25" "  -0.91%  "
Set-TextValue "E26" "  -1.39%  "
Set-TextValue "D27" "0.0000118"
Set-TextValue "E27" "  -2.11%  "
Set-TextValue "D28" "9.61"
Set-TextValue "E28" "  -6.71%  "
Set-TextValue "E29" "  +1.25%  "
Set-TextValue "D30" "0.996"
Set-TextValue "E30" "  -0.44%  "
Set-TextValue "E31" "  +1.95%  "
Set-TextValue "E32" "  -5.58%  "
Set-TextValue "D33" "1.98"
Set-TextValue "E33" "  -3.22%  "
Set-TextValue "E34" "  +0.07%  "
Set-TextValue "E35" "  -7.02%  "
Set-TextValue "D36" "7.04"
Set-TextValue "E36" "  -3.76%  "
Set-TextValue "E37" "  -2.62%  "
Set-TextValue "D38" "160.96"
Set-TextValue "E38" "  +0.27%  "
Set-TextValue "D39" "29.13"
Set-TextValue "E39" "  +11.39%  "
Set-TextValue "D40" "0.884"
Set-TextValue "E40" "  +0.93%  "
Set-TextValue "D41" "2.65"
Set-TextValue "E41" "  -4.97%  "
Set-TextValue "D42" "1.75"
Set-TextValue "E42" "  -6.19%  "
Set-TextValue "D43" "2.718.14"
Set-TextValue "E43" "  -1.34%  "
Set-TextValue "D44" "4.40"
Set-TextValue "E44" "  -5.30%  "
Set-TextValue "D45" "6.27"
Set-TextValue "E45" "  -6.81%  "
Set-TextValue "E46" "  -5.09%  "
Set-TextValue "D47" "40.24"
Set-TextValue "E47" "  -2.55%  "
Set-TextValue "D48" "0.0288"
Set-TextValue "E48" "  -3.22%  "
Set-TextValue "D49" "23.97"
Set-TextValue "E49" "  -8.26%  "
Set-TextValue "D50" "307.03"
Set-TextValue "E50" "  -5.23%  "
Set-TextValue "B51" "ONDO"
Set-TextValue "C51" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D51" "0.980"
Set-TextValue "E51" "  -5.81%  "
